$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Cells.Item(2, 4)
$cell.Value = "'50.970.26"
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -0.93%  '
$cell = $ws.Cells.Item(3, 4)
$cell.Value = "'2.929.29"
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -1.60%  '
$cell = $ws.Cells.Item(4, 4)
$cell.Value = "'0.999"
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  -0.04%  '
$cell = $ws.Cells.Item(5, 4)
$cell.Value = "'373.62"
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -1.72%  '
$cell = $ws.Cells.Item(6, 4)
$cell.Value = "'101.40"
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  -3.76%  '
$cell = $ws.Cells.Item(7, 4)
$cell.Value = "'0.534"
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -1.78%  '
$ws.Cells.Item(8, 5).Value = '  +0.01%  '
$cell = $ws.Cells.Item(9, 4)
$cell.Value = "'0.580"
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -2.64%  '
$cell = $ws.Cells.Item(10, 4)
$cell.Value = "'36.27"
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -2.62%  '
$ws.Cells.Item(11, 5).Value = '  -0.56%  '
$cell = $ws.Cells.Item(12, 4)
$cell.Value = "'0.0832"
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -1.28%  '
$cell = $ws.Cells.Item(13, 4)
$cell.Value = "'3.390.16"
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -1.90%  '
$cell = $ws.Cells.Item(14, 4)
$cell.Value = "'17.83"
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -3.96%  '
$cell = $ws.Cells.Item(15, 4)
$cell.Value = "'7.30"
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -2.42%  '
$cell = $ws.Cells.Item(16, 4)
$cell.Value = "'2.923.17"
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -2.11%  '
$cell = $ws.Cells.Item(17, 4)
$cell.Value = "'0.969"
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -0.81%  '
$cell = $ws.Cells.Item(18, 4)
$cell.Value = "'50.913.86"
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -1.04%  '
$cell = $ws.Cells.Item(19, 4)
$cell.Value = "'3.13"
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -6.97%  '
$cell = $ws.Cells.Item(20, 4)
$cell.Value = "'7.12"
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -3.67%  '
$cell = $ws.Cells.Item(21, 4)
$cell.Value = "'12.47"
$cell.Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  -4.00%  '
$cell = $ws.Cells.Item(22, 4)
$cell.Value = "'0.0₃0950"
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -0.94%  '
$cell = $ws.Cells.Item(23, 4)
$cell.Value = "'263.70"
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.61%  '
$cell = $ws.Cells.Item(24, 4)
$cell.Value = "'68.08"
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  -1.48%  '
$ws.Cells.Item(25, 5).Value = '  +2.81%  '
$cell = $ws.Cells.Item(26, 4)
$cell.Value = "'7.76"
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +7.82%  '
$cell = $ws.Cells.Item(27, 4)
$cell.Value = "'7.93"
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +6.76%  '
$cell = $ws.Cells.Item(28, 4)
$cell.Value = "'0.167"
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -1.68%  '
$ws.Cells.Item(29, 5).Value = '  -0.01%  '
$cell = $ws.Cells.Item(30, 4)
$cell.Value = "'0.112"
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -0.59%  '
$cell = $ws.Cells.Item(31, 4)
$cell.Value = "'25.56"
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -1.62%  '
$cell = $ws.Cells.Item(32, 4)
$cell.Value = "'9.80"
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.40%  '
$cell = $ws.Cells.Item(33, 4)
$cell.Value = "'50.88"
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -0.97%  '
$cell = $ws.Cells.Item(34, 4)
$cell.Value = "'0.0451"
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -1.59%  '
$cell = $ws.Cells.Item(35, 4)
$cell.Value = "'33.37"
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -4.21%  '
$ws.Cells.Item(36, 5).Value = '  -3.36%  '
$ws.Cells.Item(37, 5).Value = '  -0.20%  '
$cell = $ws.Cells.Item(38, 4)
$cell.Value = "'2.95"
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -4.44%  '
$cell = $ws.Cells.Item(39, 4)
$cell.Value = "'2.52"
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -2.95%  '
$ws.Cells.Item(40, 5).Value = '  -1.86%  '
$cell = $ws.Cells.Item(41, 4)
$cell.Value = "'16.25"
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -6.74%  '
$ws.Cells.Item(42, 5).Value = '  -4.29%  '
$cell = $ws.Cells.Item(43, 4)
$cell.Value = "'121.60"
$cell.Style = "Normal"
$ws.Cells.Item(43, 5).Value = '  -1.79%  '
$cell = $ws.Cells.Item(44, 4)
$cell.Value = "'20.79"
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -6.35%  '
$cell = $ws.Cells.Item(45, 4)
$cell.Value = "'2.04"
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -1.72%  '
$cell = $ws.Cells.Item(46, 4)
$cell.Value = "'0.271"
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -2.97%  '
$cell = $ws.Cells.Item(47, 4)
$cell.Value = "'2.29"
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -4.43%  '
$cell = $ws.Cells.Item(48, 4)
$cell.Value = "'3.19"
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -1.35%  '
$cell = $ws.Cells.Item(49, 4)
$cell.Value = "'1.979.96"
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -2.91%  '
$cell = $ws.Cells.Item(50, 4)
$cell.Value = "'0.0342"
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -4.16%  '
$cell = $ws.Cells.Item(51, 4)
$cell.Value = "'5.01"
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -3.07%  '
